$d = $word.ActiveDocument

$p19 = $d.Paragraphs(19)
$xml19 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>RF-010: Historial de Aportaciones</w:t></w:r><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>: Los usuarios deben poder ver un historial de todas sus aportaciones realizadas en cada cadena de ahorro.</w:t></w:r></w:p>'
$p19.Range.InsertXML($xml19)

$p18 = $d.Paragraphs(18)
$xml18 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>RF-009: Visualización de Cadenas Activas</w:t></w:r><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>: Los usuarios deben poder visualizar todas las cadenas de ahorro en las que están participando, con detalles como el estado actual, participantes y próximos pagos.</w:t></w:r></w:p>'
$p18.Range.InsertXML($xml18)

$p17 = $d.Paragraphs(17)
$xml17 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>RF-008: Aceptación de Invitaciones</w:t></w:r><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>: Los usuarios invitados deben poder aceptar o rechazar una invitación para unirse a una cadena de ahorro.</w:t></w:r></w:p>'
$p17.Range.InsertXML($xml17)

$p16 = $d.Paragraphs(16)
$xml16 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>RF-007: Invitación de Participantes</w:t></w:r><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>: El creador de la cadena debe poder invitar a otros usuarios a unirse a la cadena mediante correo electrónico.</w:t></w:r></w:p>'
$p16.Range.InsertXML($xml16)

$p15 = $d.Paragraphs(15)
$xml15 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>Fecha de inicio y fin de la cadena.</w:t></w:r></w:p>'
$p15.Range.InsertXML($xml15)

$p14 = $d.Paragraphs(14)
$xml14 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>Frecuencia de ahorro (semanal, quincenal, mensual).</w:t></w:r></w:p>'
$p14.Range.InsertXML($xml14)

$p13 = $d.Paragraphs(13)
$xml13 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>Monto de ahorro por participante.</w:t></w:r></w:p>'
$p13.Range.InsertXML($xml13)

$p12 = $d.Paragraphs(12)
$xml12 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>Número de participantes.</w:t></w:r></w:p>'
$p12.Range.InsertXML($xml12)

$p11 = $d.Paragraphs(11)
$xml11 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>Nombre de la cadena.</w:t></w:r></w:p>'
$p11.Range.InsertXML($xml11)

$p10 = $d.Paragraphs(10)
$xml10 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00C40C2F" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>RF-006: Creación de Cadena de Ahorro</w:t></w:r><w:r w:rsidRPr="00C40C2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>: Los usuarios deben poder crear una nueva cadena de ahorro especificando los siguientes detalles:</w:t></w:r></w:p>'
$p10.Range.InsertXML($xml10)

$p8 = $d.Paragraphs(8)
$xml8 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C40C2F" w:rsidRPr="00B3394D" w:rsidRDefault="00C40C2F" w:rsidP="00C40C2F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr></w:pPr><w:r w:rsidRPr="00B3394D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>RF-005: Perfil de Usuario</w:t></w:r><w:r w:rsidRPr="00B3394D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:eastAsia="es-CO"/></w:rPr><w:t>: Los usuarios deben poder actualizar su información personal (nombre, correo, contraseña).</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$p8.Range.InsertXML($xml8)
